# Auto-generated Excel COM-interop edit script
# Applies numeric updates to leve-profit columns (H-N) across all 8 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    @(15, 8, 1350.5807),
    @(15, 9, 1350.5807),
    @(15, 11, 4051.7421),
    @(15, 13, -3882.7421),
    @(17, 8, 1936.3914),
    @(17, 10, 1936.3914),
    @(17, 12, 5809.174199999999),
    @(17, 14, -6145.174199999999),
    @(18, 8, 712.94446),
    @(18, 9, 604.125),
    @(18, 11, 604.125),
    @(18, 13, -320.125),
    @(33, 8, 184.5238),
    @(33, 9, 144.11111),
    @(33, 11, 144.11111),
    @(33, 13, 84.88889),
    @(40, 8, 4069.353),
    @(40, 9, 3654.9092),
    @(40, 10, 4829.1665),
    @(40, 11, 3654.9092),
    @(40, 12, 4829.1665),
    @(40, 13, -3479.9092),
    @(40, 14, -5179.1665),
    @(43, 8, 1648.6),
    @(43, 9, 1066.3334),
    @(43, 10, 1898.1428),
    @(43, 11, 1066.3334),
    @(43, 12, 1898.1428),
    @(43, 13, -997.3334),
    @(43, 14, -2036.1428),
    @(54, 8, 22280.934),
    @(54, 10, 22280.934),
    @(54, 12, 22280.934),
    @(54, 14, -23252.934),
    @(59, 8, 6999.5),
    @(59, 9, 7999),
    @(59, 10, 6000),
    @(59, 11, 23997),
    @(59, 12, 18000),
    @(59, 13, -23440),
    @(59, 14, -19114),
    @(61, 8, 1684.1666),
    @(61, 9, 1684.1666),
    @(61, 11, 5052.4998),
    @(61, 13, -4880.4998),
    @(62, 8, 5395.2104),
    @(62, 9, 4994.25),
    @(62, 10, 7533.6665),
    @(62, 11, 4994.25),
    @(62, 12, 7533.6665),
    @(62, 13, -4370.25),
    @(62, 14, -8781.666499999999),
    @(64, 8, 4173.92),
    @(64, 10, 5854.125),
    @(64, 12, 5854.125),
    @(64, 14, -6350.125),
    @(65, 8, 5395.2104),
    @(65, 9, 4994.25),
    @(65, 10, 7533.6665),
    @(65, 11, 24971.25),
    @(65, 12, 37668.3325),
    @(65, 13, -21851.25),
    @(65, 14, -43908.3325),
    @(67, 8, 4173.92),
    @(67, 10, 5854.125),
    @(67, 12, 5854.125),
    @(67, 14, -7570.125),
    @(86, 8, 8499.5),
    @(86, 9, 8499.5),
    @(86, 11, 8499.5),
    @(86, 13, -7376.5),
    @(87, 8, 40000),
    @(87, 9, 0),
    @(87, 10, 40000),
    @(87, 11, 0),
    @(87, 12, 40000),
    @(87, 13, $null),
    @(87, 14, -42496),
    @(89, 8, 8499.5),
    @(89, 9, 8499.5),
    @(89, 11, 42497.5),
    @(89, 13, -36881.5),
    @(90, 8, 40000),
    @(90, 9, 0),
    @(90, 10, 40000),
    @(90, 11, 0),
    @(90, 12, 120000),
    @(90, 13, $null),
    @(90, 14, -132480),
    @(107, 8, 534.2143),
    @(107, 9, 473.25),
    @(107, 11, 473.25),
    @(107, 13, 1446.75),
    @(108, 8, 75000),
    @(108, 10, 75000),
    @(108, 12, 75000),
    @(108, 14, -82680),
    @(109, 8, 50000),
    @(109, 10, 50000),
    @(109, 12, 50000),
    @(109, 14, -52774),
    @(113, 8, 5633),
    @(113, 9, 5621.5293),
    @(113, 11, 5621.5293),
    @(113, 13, -2367.5293),
    @(129, 8, 2233.9285),
    @(129, 10, 2987.4285),
    @(129, 12, 8962.2855),
    @(129, 14, -18962.2855),
    @(132, 8, 4753.185),
    @(132, 9, 4753.185),
    @(132, 10, 0),
    @(132, 11, 14259.555),
    @(132, 12, 0),
    @(132, 13, -11729.555),
    @(132, 14, $null),
    @(134, 8, 35998.875),
    @(134, 10, 35998.875),
    @(134, 12, 35998.875),
    @(134, 14, -46138.875),
    @(135, 8, 4164),
    @(135, 10, 6000),
    @(135, 12, 54000),
    @(135, 14, -59070),
    @(137, 8, 1581.6786),
    @(137, 9, 1166.5),
    @(137, 10, 2619.625),
    @(137, 11, 3499.5),
    @(137, 12, 7858.875),
    @(137, 13, -949.5),
    @(137, 14, -12958.875),
    @(138, 8, 1953.39),
    @(138, 9, 1565.6666),
    @(138, 10, 1965.3815),
    @(138, 11, 4696.9998),
    @(138, 12, 5896.1445),
    @(138, 13, 443.0002000000004),
    @(138, 14, -16176.1445)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    @(2, 8, 1667.9),
    @(2, 9, 1525.5714),
    @(2, 10, 2000),
    @(2, 11, 1525.5714),
    @(2, 12, 2000),
    @(2, 13, -1412.5714),
    @(2, 14, -2226),
    @(5, 8, 411.4375),
    @(5, 10, 88.5),
    @(5, 12, 88.5),
    @(5, 14, -312.5),
    @(10, 8, 1000),
    @(10, 9, 1000),
    @(10, 11, 1000),
    @(10, 13, -830),
    @(11, 8, 18498.75),
    @(11, 10, 19665),
    @(11, 12, 19665),
    @(11, 14, -19953),
    @(28, 8, 26370.25),
    @(28, 9, 26370.25),
    @(28, 11, 26370.25),
    @(28, 13, -26178.25),
    @(31, 8, 1210),
    @(31, 9, 1210),
    @(31, 11, 1210),
    @(31, 13, -916),
    @(32, 8, 6504.047),
    @(32, 9, 4994.722),
    @(32, 11, 4994.722),
    @(32, 13, -4707.722),
    @(43, 8, 22932),
    @(43, 10, 22932),
    @(43, 12, 22932),
    @(43, 14, -23558),
    @(45, 8, 4594.5293),
    @(45, 9, 0),
    @(45, 10, 4594.5293),
    @(45, 11, 0),
    @(45, 12, 4594.5293),
    @(45, 13, $null),
    @(45, 14, -5348.5293),
    @(61, 8, 4136.857),
    @(61, 9, 2009.8125),
    @(61, 11, 2009.8125),
    @(61, 13, -1797.8125),
    @(63, 8, 1000),
    @(63, 9, 1000),
    @(63, 11, 1000),
    @(63, 13, -314),
    @(66, 8, 1000),
    @(66, 9, 1000),
    @(66, 11, 5000),
    @(66, 13, -1568),
    @(74, 8, 1088.9524),
    @(74, 9, 752.38464),
    @(74, 11, 752.38464),
    @(74, 13, 121.61536),
    @(77, 8, 1088.9524),
    @(77, 9, 752.38464),
    @(77, 11, 3761.9232),
    @(77, 13, 606.0767999999998),
    @(97, 8, 3604.75),
    @(97, 9, 1017.73334),
    @(97, 10, 7916.4443),
    @(97, 11, 1017.73334),
    @(97, 12, 7916.4443),
    @(97, 13, -521.73334),
    @(97, 14, -8908.444299999999),
    @(99, 8, 26370.25),
    @(99, 9, 26370.25),
    @(99, 11, 26370.25),
    @(99, 13, -23375.25),
    @(106, 8, 22500),
    @(106, 10, 22500),
    @(106, 12, 22500),
    @(106, 14, -25024),
    @(110, 8, 830.5294),
    @(110, 9, 859.3570999999999),
    @(110, 10, 696),
    @(110, 11, 859.3570999999999),
    @(110, 12, 696),
    @(110, 13, 1185.6429),
    @(110, 14, -4786),
    @(116, 8, 1667.9),
    @(116, 9, 1525.5714),
    @(116, 10, 2000),
    @(116, 11, 1525.5714),
    @(116, 12, 2000),
    @(116, 13, 768.4286),
    @(116, 14, -6588),
    @(122, 8, 4948.6113),
    @(122, 9, 5584.7036),
    @(122, 11, 16754.1108),
    @(122, 13, -14304.1108),
    @(124, 8, 184714),
    @(124, 10, 184714),
    @(124, 12, 184714),
    @(124, 14, -194534),
    @(125, 8, 139177.5),
    @(125, 10, 139177.5),
    @(125, 12, 139177.5),
    @(125, 14, -149017.5),
    @(130, 8, 37328.285),
    @(130, 10, 37328.285),
    @(130, 12, 37328.285),
    @(130, 14, -47368.285),
    @(132, 8, 4622.9473),
    @(132, 9, 4813.879),
    @(132, 11, 14441.637),
    @(132, 13, -11911.637),
    @(136, 8, 4136.857),
    @(136, 9, 2009.8125),
    @(136, 11, 6029.4375),
    @(136, 13, -3479.4375)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    @(3, 8, 1667.9),
    @(3, 9, 1525.5714),
    @(3, 10, 2000),
    @(3, 11, 1525.5714),
    @(3, 12, 2000),
    @(3, 13, -1411.5714),
    @(3, 14, -2228),
    @(4, 8, 411.4375),
    @(4, 10, 88.5),
    @(4, 12, 88.5),
    @(4, 14, -318.5),
    @(19, 8, 7000),
    @(19, 10, 7000),
    @(19, 12, 7000),
    @(19, 14, -7346),
    @(20, 8, 1856.5454),
    @(20, 9, 2135.375),
    @(20, 10, 1113),
    @(20, 11, 2135.375),
    @(20, 12, 1113),
    @(20, 13, -1888.375),
    @(20, 14, -1607),
    @(26, 8, 33666.332),
    @(26, 9, 33666.332),
    @(26, 11, 33666.332),
    @(26, 13, -33374.332),
    @(94, 8, 650.1905),
    @(94, 9, 300.85715),
    @(94, 10, 1348.8572),
    @(94, 11, 300.85715),
    @(94, 12, 1348.8572),
    @(94, 13, 150.14285),
    @(94, 14, -2250.8572),
    @(96, 8, 22441.5),
    @(96, 9, 18888),
    @(96, 11, 18888),
    @(96, 13, -16142),
    @(107, 8, 1813.5555),
    @(107, 9, 1902),
    @(107, 11, 1902),
    @(107, 13, 18),
    @(132, 8, 80999.5),
    @(132, 10, 80999.5),
    @(132, 12, 80999.5),
    @(132, 14, -91119.5),
    @(134, 8, 6534.245),
    @(134, 9, 4959.6855),
    @(134, 11, 14879.0565),
    @(134, 13, -12344.0565)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    @(16, 8, 216.5),
    @(16, 9, 216.5),
    @(16, 10, 0),
    @(16, 11, 216.5),
    @(16, 12, 0),
    @(16, 13, 70.5),
    @(16, 14, $null),
    @(31, 8, 43158.84),
    @(31, 9, 1745.7333),
    @(31, 10, 105278.5),
    @(31, 11, 1745.7333),
    @(31, 12, 105278.5),
    @(31, 13, -1450.7333),
    @(31, 14, -105868.5),
    @(34, 8, 43158.84),
    @(34, 9, 1745.7333),
    @(34, 10, 105278.5),
    @(34, 11, 1745.7333),
    @(34, 12, 105278.5),
    @(34, 13, -1543.7333),
    @(34, 14, -105682.5),
    @(51, 8, 19117.47),
    @(51, 10, 49998.5),
    @(51, 12, 49998.5),
    @(51, 14, -51470.5),
    @(61, 8, 19117.47),
    @(61, 10, 49998.5),
    @(61, 12, 49998.5),
    @(61, 14, -50694.5),
    @(74, 8, 39999.75),
    @(74, 10, 39999.75),
    @(74, 12, 39999.75),
    @(74, 14, -41747.75),
    @(77, 8, 39999.75),
    @(77, 10, 39999.75),
    @(77, 12, 119999.25),
    @(77, 14, -128735.25),
    @(86, 8, 6010.294),
    @(86, 9, 5808.222),
    @(86, 10, 6237.625),
    @(86, 11, 5808.222),
    @(86, 12, 6237.625),
    @(86, 13, -4685.222),
    @(86, 14, -8483.625),
    @(89, 8, 6010.294),
    @(89, 9, 5808.222),
    @(89, 10, 6237.625),
    @(89, 11, 29041.11),
    @(89, 12, 31188.125),
    @(89, 13, -23425.11),
    @(89, 14, -42420.125),
    @(107, 8, 1100),
    @(107, 9, 1000),
    @(107, 11, 1000),
    @(107, 13, 920),
    @(113, 8, 216.5),
    @(113, 9, 216.5),
    @(113, 10, 0),
    @(113, 11, 216.5),
    @(113, 12, 0),
    @(113, 13, 1953.5),
    @(113, 14, $null),
    @(132, 8, 3781.2058),
    @(132, 9, 3549.3809),
    @(132, 10, 4155.6924),
    @(132, 11, 10648.1427),
    @(132, 12, 12467.0772),
    @(132, 13, -8118.1427),
    @(132, 14, -17527.0772),
    @(133, 8, 93653.8),
    @(133, 10, 94567.25),
    @(133, 12, 94567.25),
    @(133, 14, -99627.25),
    @(134, 8, 102969.6),
    @(134, 9, 125899.5),
    @(134, 11, 377698.5),
    @(134, 13, -375163.5),
    @(141, 8, 262793.28),
    @(141, 10, 445621),
    @(141, 12, 445621),
    @(141, 14, -455981)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$edits = @(
    @(4, 8, 18634880),
    @(4, 9, 22161856),
    @(4, 11, 66485568),
    @(4, 13, -66485456),
    @(7, 8, 204.2),
    @(7, 9, 114.57143),
    @(7, 11, 343.71429),
    @(7, 13, -231.71429),
    @(8, 8, 951),
    @(8, 9, 951),
    @(8, 11, 2853),
    @(8, 13, -2714),
    @(23, 8, 285.81818),
    @(23, 9, 135.25),
    @(23, 11, 405.75),
    @(23, 13, -170.75),
    @(37, 8, 139601.95),
    @(37, 10, 139601.95),
    @(37, 12, 418805.85),
    @(37, 14, -419029.85),
    @(107, 8, 1668.7142),
    @(107, 9, 665.5),
    @(107, 10, 2070),
    @(107, 11, 1996.5),
    @(107, 12, 6210),
    @(107, 13, -76.5),
    @(107, 14, -10050),
    @(117, 8, 4232.8),
    @(117, 10, 4232.8),
    @(117, 12, 12698.4),
    @(117, 14, -19582.4),
    @(129, 8, 2779.4075),
    @(129, 9, 549.8333),
    @(129, 10, 3416.4285),
    @(129, 11, 1649.4999),
    @(129, 12, 10249.2855),
    @(129, 13, 3350.5001),
    @(129, 14, -20249.2855),
    @(131, 8, 1593.0714),
    @(131, 9, 1146.25),
    @(131, 10, 2188.8333),
    @(131, 11, 3438.75),
    @(131, 12, 6566.499899999999),
    @(131, 13, 1601.25),
    @(131, 14, -16646.4999),
    @(132, 8, 4935.6665),
    @(132, 9, 2330.125),
    @(132, 10, 7913.4287),
    @(132, 11, 20971.125),
    @(132, 12, 71220.85830000001),
    @(132, 13, -18441.125),
    @(132, 14, -76280.85830000001)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$edits = @(
    @(2, 8, 453),
    @(2, 9, 65),
    @(2, 10, 986.5),
    @(2, 11, 65),
    @(2, 12, 986.5),
    @(2, 13, 48),
    @(2, 14, -1212.5),
    @(26, 8, 39495),
    @(26, 9, 39000),
    @(26, 10, 39990),
    @(26, 11, 39000),
    @(26, 12, 39990),
    @(26, 13, -38720),
    @(26, 14, -40550),
    @(50, 8, 39495),
    @(50, 9, 39000),
    @(50, 10, 39990),
    @(50, 11, 39000),
    @(50, 12, 39990),
    @(50, 13, -38502),
    @(50, 14, -40986),
    @(63, 8, 0),
    @(63, 9, 0),
    @(63, 10, 0),
    @(63, 11, 0),
    @(63, 12, 0),
    @(63, 13, $null),
    @(63, 14, $null),
    @(66, 8, 0),
    @(66, 9, 0),
    @(66, 10, 0),
    @(66, 11, 0),
    @(66, 12, 0),
    @(66, 13, $null),
    @(66, 14, $null),
    @(70, 8, 10424.923),
    @(70, 9, 11716.25),
    @(70, 11, 11716.25),
    @(70, 13, -11446.25),
    @(73, 8, 10424.923),
    @(73, 9, 11716.25),
    @(73, 11, 11716.25),
    @(73, 13, -10780.25),
    @(111, 8, 77570.14),
    @(111, 10, 108249.75),
    @(111, 12, 108249.75),
    @(111, 14, -114383.75),
    @(122, 8, 39359.3),
    @(122, 9, 68115.94),
    @(122, 10, 6494.5713),
    @(122, 11, 204347.82),
    @(122, 12, 19483.7139),
    @(122, 13, -201897.82),
    @(122, 14, -24383.7139),
    @(132, 8, 8164.5835),
    @(132, 9, 4899.5),
    @(132, 10, 8817.6),
    @(132, 11, 14698.5),
    @(132, 12, 26452.8),
    @(132, 13, -12168.5),
    @(132, 14, -31512.8)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    @(7, 8, 8519.429),
    @(7, 9, 9481),
    @(7, 11, 9481),
    @(7, 13, -9369),
    @(16, 8, 980.4211),
    @(16, 9, 929.3333),
    @(16, 11, 929.3333),
    @(16, 13, -759.3333),
    @(22, 8, 901.9091),
    @(22, 9, 443.83334),
    @(22, 11, 443.83334),
    @(22, 13, -148.83334),
    @(25, 8, 1999),
    @(25, 9, 0),
    @(25, 10, 1999),
    @(25, 11, 0),
    @(25, 12, 1999),
    @(25, 13, $null),
    @(25, 14, -2459),
    @(27, 8, 901.9091),
    @(27, 9, 443.83334),
    @(27, 11, 443.83334),
    @(27, 13, -336.83334),
    @(40, 8, 3266),
    @(40, 9, 2749),
    @(40, 11, 2749),
    @(40, 13, -2613),
    @(61, 8, 2164),
    @(61, 9, 2164),
    @(61, 11, 2164),
    @(61, 13, -1962),
    @(68, 8, 2531.5334),
    @(68, 9, 2547.3),
    @(68, 11, 2547.3),
    @(68, 13, -1798.3),
    @(71, 8, 2531.5334),
    @(71, 9, 2547.3),
    @(71, 11, 12736.5),
    @(71, 13, -8992.5),
    @(82, 8, 3517.5),
    @(82, 9, 1227.125),
    @(82, 11, 1227.125),
    @(82, 13, -866.125),
    @(85, 8, 3517.5),
    @(85, 9, 1227.125),
    @(85, 11, 1227.125),
    @(85, 13, 20.875),
    @(93, 8, 1731.4),
    @(93, 9, 1506),
    @(93, 11, 1506),
    @(93, 13, -258),
    @(95, 8, 32999.5),
    @(95, 10, 32999.5),
    @(95, 12, 32999.5),
    @(95, 14, -38491.5),
    @(96, 8, 48678.6),
    @(96, 10, 48678.6),
    @(96, 12, 48678.6),
    @(96, 14, -54170.6),
    @(113, 8, 2164),
    @(113, 9, 2164),
    @(113, 11, 2164),
    @(113, 13, 6),
    @(126, 8, 8519.429),
    @(126, 9, 9481),
    @(126, 11, 28443),
    @(126, 13, -25973),
    @(132, 8, 4086.8684),
    @(132, 9, 3787.652),
    @(132, 11, 11362.956),
    @(132, 13, -8832.956),
    @(136, 8, 3769.8125),
    @(136, 9, 3909.7368),
    @(136, 11, 11729.2104),
    @(136, 13, -9179.2104)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    @(3, 8, 4504752),
    @(3, 9, 9500),
    @(3, 11, 9500),
    @(3, 13, -9386),
    @(46, 8, 52801),
    @(46, 10, 52801),
    @(46, 12, 52801),
    @(46, 14, -53263),
    @(81, 8, 6020.2915),
    @(81, 9, 8399),
    @(81, 11, 16798),
    @(81, 13, -15737),
    @(84, 8, 6020.2915),
    @(84, 9, 8399),
    @(84, 11, 83990),
    @(84, 13, -78686),
    @(93, 8, 75000),
    @(93, 10, 75000),
    @(93, 12, 75000),
    @(93, 14, -79992),
    @(97, 8, 54443),
    @(97, 10, 54443),
    @(97, 12, 54443),
    @(97, 14, -56425),
    @(99, 8, 38379),
    @(99, 9, 32449.5),
    @(99, 10, 50238),
    @(99, 11, 32449.5),
    @(99, 12, 50238),
    @(99, 13, -29454.5),
    @(99, 14, -56228),
    @(107, 8, 781.2857),
    @(107, 9, 673.8),
    @(107, 10, 879),
    @(107, 11, 2021.4),
    @(107, 12, 2637),
    @(107, 13, -101.3999999999999),
    @(107, 14, -6477),
    @(113, 8, 241.18182),
    @(113, 9, 247.8),
    @(113, 11, 743.4000000000001),
    @(113, 13, 1426.6),
    @(122, 8, 89008),
    @(122, 9, 115350.4),
    @(122, 10, 1200),
    @(122, 11, 346051.2),
    @(122, 12, 3600),
    @(122, 13, -343601.2),
    @(122, 14, -8500),
    @(126, 8, 2182),
    @(126, 9, 1860.1666),
    @(126, 10, 2825.6667),
    @(126, 11, 5580.4998),
    @(126, 12, 8477.000100000001),
    @(126, 13, -3110.4998),
    @(126, 14, -13417.0001),
    @(127, 8, 72397),
    @(127, 9, 58492.5),
    @(127, 11, 58492.5),
    @(127, 13, -53532.5),
    @(128, 8, 32979.285),
    @(128, 10, 32979.285),
    @(128, 12, 32979.285),
    @(128, 14, -42939.285),
    @(132, 8, 2135.2),
    @(132, 9, 2135.2),
    @(132, 11, 6405.599999999999),
    @(132, 13, -3875.599999999999),
    @(134, 8, 52801),
    @(134, 10, 52801),
    @(134, 12, 158403),
    @(134, 14, -163473),
    @(136, 8, 6521.775),
    @(136, 9, 4933.943),
    @(136, 10, 17636.6),
    @(136, 11, 14801.829),
    @(136, 12, 52909.8),
    @(136, 13, -12251.829),
    @(136, 14, -58009.8),
    @(138, 8, 98999.5),
    @(138, 10, 98999.5),
    @(138, 12, 98999.5),
    @(138, 14, -109279.5)
)
foreach ($e in $edits) {
    $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
}

Write-Output "Applied $(654) cell edits across 8 sheets"